# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (detail holdings) right before the
#    "总计" (total) summary sheet, matching the layout of the other
#    quarterly sheets (e.g. "2021-Q4").
# 2. Add a new leading row to the "总计" sheet summarising the new quarter
#    and renumber the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet, inserted just before "总计"
# ---------------------------------------------------------------------
$insertBefore = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($insertBefore)
$newSheet.Name = "2022-Q1"

# NOTE: sheet references are tracked by position, not identity, so any
# variable captured before the Add() above (e.g. a cached "总计" handle)
# would now silently point at this new sheet instead. Re-resolve sheets
# we still need by name *after* the structural change, right before use.

# Columns B:G hold numeric-looking text (fund codes, money amounts as
# strings) in every quarterly sheet -- force them to Text *before* typing
# the values so COM doesn't silently coerce them to numbers.
$newSheet.Range("B2:G4").NumberFormat = "@"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "159869"
$newSheet.Range("C2").Value = "华夏中证动漫游戏ETF"
$newSheet.Range("D2").Value = "6.20"
$newSheet.Range("E2").Value = "98.75"
$newSheet.Range("F2").Value = "3.55"
$newSheet.Range("G2").Value = "0.2201"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "516010"
$newSheet.Range("C3").Value = "国泰中证动漫游戏ETF"
$newSheet.Range("D3").Value = "4.95"
$newSheet.Range("E3").Value = "98.91"
$newSheet.Range("F3").Value = "3.49"
$newSheet.Range("G3").Value = "0.1728"
$newSheet.Range("H3").Value = 9

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "516770"
$newSheet.Range("C4").Value = "华泰柏瑞中证动漫游戏ETF"
$newSheet.Range("D4").Value = "1.11"
$newSheet.Range("E4").Value = "96.56"
$newSheet.Range("F4").Value = "3.50"
$newSheet.Range("G4").Value = "0.0388"
$newSheet.Range("H4").Value = 9

# Pick up the header/index-column look (bold, centered, bordered) from the
# sibling "2021-Q4" sheet without disturbing the values/number formats we
# just set.
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Range("A1:H4").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

# Re-apply the index-column style (bold/centered/bordered) that column A
# carries on every other data row.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.43

# Renumber the zero-based index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

$totalSheet.Range("A1").Select()
